$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 4.9372844376893701
$ws.Range("C2").Value = 2.0588463868717524

$ws.Range("B3").Value = 1.1558967911587537
$ws.Range("C3").Value = 3.2852655384412373

$ws.Range("B4").Value = 4.3201072330525276
$ws.Range("C4").Value = 5.2624351432080676

$ws.Range("B5").Value = 2.8874878263154731
$ws.Range("C5").Value = 3.0522904438843694

$ws.Range("C6").Value = 5.5522896573485676

$ws.Range("B9").Value = 109.43629248616479
$ws.Range("C9").Value = 2.3418180302461367

$ws.Range("C10").Value = 154.60363398197694

$ws.Range("B12").Value = 36.684652282188203
$ws.Range("C12").Value = 32.237774215338959

$ws.Range("B15").Value = 0.59834630534917566
$ws.Range("C15").Value = 0.35330435208717326

$ws.Range("B16").Value = 10.012651951040695
$ws.Range("C16").Value = 8.0539692291398843
